$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4927675928165911
$ws.Range("D2").Value = 0.009140987240385812
$ws.Range("E2").Value = 0.09638581504330546
$ws.Range("F2").Value = 0.3212188075540041
$ws.Range("G2").Value = 0.2049276593427862
$ws.Range("H2").Value = 0.3008942488266513
$ws.Range("I2").Value = 1.132485420188019
$ws.Range("M2").Value = 1.432321805107307
$ws.Range("O2").Value = 0.9182964999054093
# Row 3
$ws.Range("B3").Value = 0.4329642906879485
$ws.Range("D3").Value = 0.008086000681259975
$ws.Range("E3").Value = 0.09786825833404933
$ws.Range("F3").Value = 0.3054797143202919
$ws.Range("G3").Value = 0.1893833530526194
$ws.Range("H3").Value = 0.2977200459642262
$ws.Range("I3").Value = 1.026373579380177
$ws.Range("M3").Value = 1.252898307377208
$ws.Range("O3").Value = 0.8783809609745674
# Row 4
$ws.Range("B4").Value = 0.3960877178337228
$ws.Range("D4").Value = 0.007436720435912036
$ws.Range("E4").Value = 0.09906192152865145
$ws.Range("F4").Value = 0.296108005691238
$ws.Range("G4").Value = 0.1800417492896429
$ws.Range("H4").Value = 0.2960092709603401
$ws.Range("I4").Value = 0.9615371276990601
$ws.Range("M4").Value = 1.142240772238551
$ws.Range("O4").Value = 0.8547894254240873
# Row 5
$ws.Range("B5").Value = 0.3810214938868057
$ws.Range("D5").Value = 0.007171783524981379
$ws.Range("E5").Value = 0.09961920259038948
$ws.Range("F5").Value = 0.2923619338133534
$ws.Range("G5").Value = 0.1762853024405899
$ws.Range("H5").Value = 0.2953718103503178
$ws.Range("I5").Value = 0.9351975914198078
$ws.Range("M5").Value = 1.097027049475543
$ws.Range("O5").Value = 0.8454044460618491
# Row 6
$ws.Range("B6").Value = 0.3785174438621652
$ws.Range("D6").Value = 0.007127770726196303
$ws.Range("E6").Value = 0.09971600594746022
$ws.Range("F6").Value = 0.2917442946282804
$ws.Range("G6").Value = 0.1756645688524827
$ws.Range("H6").Value = 0.2952695599473572
$ws.Range("I6").Value = 0.9308289417169817
$ws.Range("M6").Value = 1.089512210258718
$ws.Range("O6").Value = 0.8438598424241661
# Row 7
$ws.Range("B7").Value = 0.3958846851888609
$ws.Range("D7").Value = 0.007433148779128373
$ws.Range("E7").Value = 0.09906915088799195
$ws.Range("F7").Value = 0.2960571900315969
$ws.Range("G7").Value = 0.1799908855723231
$ws.Range("H7").Value = 0.296000432515811
$ws.Range("I7").Value = 0.9611815691098116
$ws.Range("M7").Value = 1.141631484803099
$ws.Range("O7").Value = 0.8546619324019957
# Row 8
$ws.Range("B8").Value = 0.4721804648393686
$ws.Range("D8").Value = 0.008777560289772168
$ws.Range("E8").Value = 0.09683787926031329
$ws.Range("F8").Value = 0.3157310306950336
$ws.Range("G8").Value = 0.1995255370700448
$ws.Range("H8").Value = 0.2997502216866366
$ws.Range("I8").Value = 1.095833607232933
$ws.Range("M8").Value = 1.370560487283441
$ws.Range("O8").Value = 0.9043422715638485
# Row 9
$ws.Range("B9").Value = 0.6205236852527776
$ws.Range("D9").Value = 0.01140070335514309
$ws.Range("E9").Value = 0.09473082659828336
$ws.Range("F9").Value = 0.3566538820269471
$ws.Range("G9").Value = 0.2394712311944289
$ws.Range("H9").Value = 0.3090035230383279
$ws.Range("I9").Value = 1.36231893525337
$ws.Range("M9").Value = 1.81546455500802
$ws.Range("O9").Value = 1.009124801682248
# Row 10
$ws.Range("B10").Value = 0.728713051287798
$ws.Range("D10").Value = 0.01331838475246627
$ws.Range("E10").Value = 0.09459361794134935
$ws.Range("F10").Value = 0.3881849351860609
$ws.Range("G10").Value = 0.26986270244997
$ws.Range("H10").Value = 0.3169752099586987
$ws.Range("I10").Value = 1.559509024266362
$ws.Range("M10").Value = 2.13974165898
$ws.Range("O10").Value = 1.090721609134789
# Row 11
$ws.Range("B11").Value = 0.7777542647986593
$ws.Range("D11").Value = 0.0141884168868458
$ws.Range("E11").Value = 0.09484376353069024
$ws.Range("F11").Value = 0.4028554748166044
$ws.Range("G11").Value = 0.2839248622061916
$ws.Range("H11").Value = 0.3208598096815365
$ws.Range("I11").Value = 1.64950678276054
$ws.Range("M11").Value = 2.2866741536929
$ws.Range("O11").Value = 1.128871459512339
# Row 12
$ws.Range("B12").Value = 0.7962992792301407
$ws.Range("D12").Value = 0.01451751139968138
$ws.Range("E12").Value = 0.09498394608602467
$ws.Range("F12").Value = 0.4084584170940104
$ws.Range("G12").Value = 0.289284627349403
$ws.Range("H12").Value = 0.3223681841704291
$ws.Range("I12").Value = 1.683627521378412
$ws.Range("M12").Value = 2.342227067881993
$ws.Range("O12").Value = 1.14346812575252
# Row 13
$ws.Range("B13").Value = 0.7923064383174392
$ws.Range("D13").Value = 0.01444665188302707
$ws.Range("E13").Value = 0.09495172525911499
$ws.Range("F13").Value = 0.4072496011825848
$ws.Range("G13").Value = 0.2881287508698023
$ws.Range("H13").Value = 0.3220416634931809
$ws.Range("I13").Value = 1.676277235728691
$ws.Range("M13").Value = 2.330266687285359
$ws.Range("O13").Value = 1.140317764779297
# Row 14
$ws.Range("B14").Value = 0.7792804959148043
$ws.Range("D14").Value = 0.01421549924740617
$ws.Range("E14").Value = 0.09485438242487731
$ws.Range("F14").Value = 0.4033154767631686
$ws.Range("G14").Value = 0.2843651134790406
$ws.Range("H14").Value = 0.3209831544935895
$ws.Range("I14").Value = 1.652313114411839
$ws.Range("M14").Value = 2.291246293005401
$ws.Range("O14").Value = 1.130069316107807
# Row 15
$ws.Range("B15").Value = 0.7712983457921609
$ws.Range("D15").Value = 0.01407386272321531
$ws.Range("E15").Value = 0.09480069243736011
$ws.Range("F15").Value = 0.4009119168020163
$ws.Range("G15").Value = 0.2820643194267092
$ws.Range("H15").Value = 0.3203396593794139
$ws.Range("I15").Value = 1.637639623505834
$ws.Range("M15").Value = 2.267333714152187
$ws.Range("O15").Value = 1.123811454629049
# Row 16
$ws.Range("B16").Value = 0.7255045163500995
$ws.Range("D16").Value = 0.01326147641405129
$ws.Range("E16").Value = 0.09458360341639604
$ws.Range("F16").Value = 0.3872328034498196
$ws.Range("G16").Value = 0.2689485313037068
$ws.Range("H16").Value = 0.3167265580134
$ws.Range("I16").Value = 1.553633232460214
$ws.Range("M16").Value = 2.130127252816436
$ws.Range("O16").Value = 1.088249334996021
# Row 17
$ws.Range("B17").Value = 0.6973662323703707
$ws.Range("D17").Value = 0.01276248351311438
$ws.Range("E17").Value = 0.09453084211077112
$ws.Range("F17").Value = 0.3789252081078445
$ws.Range("G17").Value = 0.2609635730742497
$ws.Range("H17").Value = 0.3145763396631764
$ws.Range("I17").Value = 1.502172222395984
$ws.Range("M17").Value = 2.045803802240329
$ws.Range("O17").Value = 1.066698543166098
# Row 18
$ws.Range("B18").Value = 0.6811654353848837
$ws.Range("D18").Value = 0.01247525855259823
$ws.Range("E18").Value = 0.09452988384138195
$ws.Range("F18").Value = 0.3741776379180095
$ws.Range("G18").Value = 0.2563930975486954
$ws.Range("H18").Value = 0.3133638840985498
$ws.Range("I18").Value = 1.4726010836809
$ws.Range("M18").Value = 1.997248529007521
$ws.Range("O18").Value = 1.054400032916931
# Row 19
$ws.Range("B19").Value = 0.6756773292764819
$ws.Range("D19").Value = 0.01237797276687047
$ws.Range("E19").Value = 0.09453459148851806
$ws.Range("F19").Value = 0.3725754558495282
$ws.Range("G19").Value = 0.2548494163182653
$ws.Range("H19").Value = 0.3129575329596577
$ws.Range("I19").Value = 1.462593650217997
$ws.Range("M19").Value = 1.980799266561235
$ws.Range("O19").Value = 1.050252560032789
# Row 20
$ws.Range("B20").Value = 0.7003633033467622
$ws.Range("D20").Value = 0.01281562489997157
$ws.Range("E20").Value = 0.09453341334691601
$ws.Range("F20").Value = 0.3798063806908232
$ws.Range("G20").Value = 0.261811276766295
$ws.Range("H20").Value = 0.3148027181851631
$ws.Range("I20").Value = 1.507647463902515
$ws.Range("M20").Value = 2.054785859166714
$ws.Range("O20").Value = 1.068982615680312
# Row 21
$ws.Range("B21").Value = 0.7831072374628434
$ws.Range("D21").Value = 0.01428340462066302
$ws.Range("E21").Value = 0.09488173653432597
$ws.Range("F21").Value = 0.4044697304617415
$ws.Range("G21").Value = 0.2854696371772576
$ws.Range("H21").Value = 0.3212930484766758
$ws.Range("I21").Value = 1.659350872011061
$ws.Range("M21").Value = 2.302709923073422
$ws.Range("O21").Value = 1.133075444002174
# Row 22
$ws.Range("B22").Value = 0.8370343220594805
$ws.Range("D22").Value = 0.01524052481222782
$ws.Range("E22").Value = 0.09537459352418409
$ws.Range("F22").Value = 0.4208659210494119
$ws.Range("G22").Value = 0.3011344850086459
$ws.Range("H22").Value = 0.3257527275306131
$ws.Range("I22").Value = 1.758733312389552
$ws.Range("M22").Value = 2.464232637160535
$ws.Range("O22").Value = 1.175839849939422
# Row 23
$ws.Range("B23").Value = 0.8082664532821582
$ws.Range("D23").Value = 0.01472990001360586
$ws.Range("E23").Value = 0.09508710831526912
$ws.Range("F23").Value = 0.4120894245148605
$ws.Range("G23").Value = 0.292755091926395
$ws.Range("H23").Value = 0.323352502642976
$ws.Range("I23").Value = 1.705670083641394
$ws.Range("M23").Value = 2.378072701442619
$ws.Range("O23").Value = 1.152934902711081
# Row 24
$ws.Range("B24").Value = 0.6990084014010449
$ws.Range("D24").Value = 0.01279160075859664
$ws.Range("E24").Value = 0.09453215941884707
$ws.Range("F24").Value = 0.3794079135898443
$ws.Range("G24").Value = 0.2614279671098956
$ws.Range("H24").Value = 0.3147002985402736
$ws.Range("I24").Value = 1.505172061967301
$ws.Range("M24").Value = 2.050725309586568
$ws.Range("O24").Value = 1.067949702159581
# Row 25
$ws.Range("B25").Value = 0.5805316023148919
$ws.Range("D25").Value = 0.01069264379918877
$ws.Range("E25").Value = 0.09505513364047502
$ws.Range("F25").Value = 0.3453285596630593
$ws.Range("G25").Value = 0.2284848564421651
$ws.Range("H25").Value = 0.3062953034852569
$ws.Range("I25").Value = 1.289975735151188
$ws.Range("M25").Value = 1.69555192218715
$ws.Range("O25").Value = 0.979977380432814

Write-Host "Updated pl_mw values for 380 kV case"
